$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42: B42 switches from a text "3" to a real number 3 (everything else in the row is unchanged)
$ws.Range("B42").Value = 3

# New row 43: a fresh annotation row appended right after row 42
$ws.Range("A43").Value = "Ying Tang"

# Keep B43 as a text "3" (matching the format the old B42 used before this edit),
# forcing text storage then restoring the default style so no stray number format sticks.
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "3"
$ws.Range("B43").Style = "Normal"

$ws.Range("C43").Value = "无"
$ws.Range("D43").Value = "DFT"
$ws.Range("E43").Value = "MET"
$ws.Range("F43").Value = "cf97de89-8b46-4ca2-a071-801296a106cf"
$ws.Range("G43").Value = "SkYXvCR6W_annotated.xlsx"
$ws.Range("H43").Value = "Some important implementation details are missing (activation functions, loss function used), and others have to be deduced by observing the output dimensions of the individual layers of the network."
